# Edit script: renames several header labels and updates column C (GDP) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text (shared-string) renames ---
# Map of old header text -> new header text. We look these up by scanning row 1
# so the change is applied regardless of exact column position.
$headerRenames = @{
    "Gross_National_Income"                 = "GDP"
    "NGO_Country_Budget_Previous_Year"       = "Budget_Previous_Year"
    "Vision_ONGD_LatinAmerica"               = "LatinAmerica"
    "Vision_ONGD_Africa"                     = "Africa"
    "Vision_ONGD_Confessional"               = "Confessional"
    "Vision_ONGD_Universal"                  = "Universal"
    "Total_subvencion_en_el_Pais_y_Anyo"     = "Donor_Aid_Budget"
    "Total_Fondos"                           = "Total_Funds"
    "Proporcion_Fondos_Privados"             = "%_Private_Funds"
    "Proporcion_Fondos_MAE"                  = "%_MAE_Funds"
    "Delegacion"                             = "Delegation"
}

$usedRange = $ws.UsedRange
$headerRow = 1
$lastCol = $usedRange.Columns.Count

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item($headerRow, $c)
    $val = $cell.Value2
    if ($headerRenames.ContainsKey($val)) {
        $cell.Value = $headerRenames[$val]
    }
}

# --- Column C (GDP) numeric updates for rows 2..46 ---
$gdpValues = @(
    5596.139681459835,
    2934.187009790061,
    2870.311589353206,
    1460.056109840828,
    4729.735976516416,
    1909.084588129339,
    12227.21453003286,
    6128.19547247793,
    4547.50930098406,
    3972.630273980753,
    2983.242707849043,
    2898.942214704482,
    1503.870423231357,
    1955.461557360978,
    6336.709213679884,
    4355.934938677345,
    5082.354756663512,
    3083.80337578809,
    2965.153206179127,
    1577.487171555845,
    2024.117324382548,
    6711.616186806423,
    4479.398934239905,
    5360.226632400601,
    2999.422762626143,
    4394.543881413723,
    5642.578115155247,
    3056.152683606517,
    4699.493713911862,
    5919.20956823756,
    5412.131646018807,
    3252.634165082374,
    3137.260298393558,
    4861.287098802361,
    5996.49696468919,
    5330.539154475424,
    3314.741082534716,
    3210.869677115934,
    4944.191641077407,
    6114.227214287786,
    5176.058803160127,
    3382.563653843273,
    3242.636921959078,
    5089.61202008711,
    6262.368904654469
)

$startRow = 2
for ($i = 0; $i -lt $gdpValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = $gdpValues[$i]
}
